# Daily update at 8 AM UTC
# Appends the next day's row (day 90) to the "Wins Over Time" tracking sheet,
# continuing directly after the existing last row (row 89 / 2026-01-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count   # 89
$newRow  = $lastRow + 1               # 90

# New day's date is one day after the previous row's date.
$ws.Range("A$newRow").Value = $ws.Range("A$lastRow").Value2 + 1
$ws.Range("B$newRow").Value = 210
$ws.Range("C$newRow").Value = 215
$ws.Range("D$newRow").Value = 201

# Match the date-column formatting used by the rest of the table.
$ws.Range("A$newRow").NumberFormat = $ws.Range("A$lastRow").NumberFormat
